# Adds two new survey-wave columns (13. 4. 2021 and 4. 5. 2021) to both
# worksheets, and refreshes the "aktualizace" date in the two footer notes.
$wb = $excel.ActiveWorkbook

# --- Sheet "data": add columns AB/AC for the two new survey waves ---
$ws1 = $wb.Worksheets.Item("data")
$ws1.Range("AB1").Value = "13. 4. 2021"
$ws1.Range("AC1").Value = "4. 5. 2021"
$ws1.Range("AA1").Copy()
$ws1.Range("AB1:AC1").PasteSpecial(-4122)

$ws1.Range("AB2").Value = 0.2
$ws1.Range("AC2").Value = 0.2
$ws1.Range("AB3").Value = 0.12
$ws1.Range("AC3").Value = 0.11
$ws1.Range("AB4").Value = 0.48
$ws1.Range("AC4").Value = 0.53
$ws1.Range("AB5").Value = 0.28
$ws1.Range("AC5").Value = 0.27
$ws1.Range("AB6").Value = 0.13
$ws1.Range("AC6").Value = 0.12
$ws1.Range("AB7").Value = 0.2
$ws1.Range("AC7").Value = 0.19
$ws1.Range("AB8").Value = 0.3
$ws1.Range("AC8").Value = 0.25
$ws1.Range("AB9").Value = 0.16
$ws1.Range("AC9").Value = 0.17
$ws1.Range("AB10").Value = 0.23
$ws1.Range("AC10").Value = 0.22
$ws1.Range("AB11").Value = 0.19
$ws1.Range("AC11").Value = 0.19
$ws1.Range("AB12").Value = 0.26
$ws1.Range("AC12").Value = 0.21
$ws1.Range("AB13").Value = 0.41
$ws1.Range("AC13").Value = 0.35
$ws1.Range("AB14").Value = 0.17
$ws1.Range("AC14").Value = 0.18
$ws1.Range("AB15").Value = 0.2
$ws1.Range("AC15").Value = 0.17
$ws1.Range("AB16").Value = 0.23
$ws1.Range("AC16").Value = 0.23
$ws1.Range("AB17").Value = 0.18
$ws1.Range("AC17").Value = 0.19
$ws1.Range("AB18").Value = 0.23
$ws1.Range("AC18").Value = 0.26
$ws1.Range("AB19").Value = 0.27
$ws1.Range("AC19").Value = 0.24
$ws1.Range("AB20").Value = 0.15
$ws1.Range("AC20").Value = 0.15
$ws1.Range("AB21").Value = 0.12
$ws1.Range("AC21").Value = 0.14
$ws1.Range("AB22").Value = 0.11
$ws1.Range("AC22").Value = 0.13
$ws1.Range("AB23").Value = 0.21
$ws1.Range("AC23").Value = 0.19
$ws1.Range("AB24").Value = 0.41
$ws1.Range("AC24").Value = 0.41
$ws1.Range("AB25").Value = 0.36
$ws1.Range("AC25").Value = 0.41
$ws1.Range("AB26").Value = 0.16
$ws1.Range("AC26").Value = 0.13
$ws1.Range("AB27").Value = 0.06
$ws1.Range("AC27").Value = 0.06
$ws1.Range("AB28").Value = 0.12
$ws1.Range("AC28").Value = 0.12
$ws1.Range("AB29").Value = 0.24
$ws1.Range("AC29").Value = 0.15
$ws1.Range("AB30").Value = 0.08
$ws1.Range("AC30").Value = 0.07
$ws1.Range("AB31").Value = 0.11
$ws1.Range("AC31").Value = 0.1
$ws1.Range("AB32").Value = 0.12
$ws1.Range("AC32").Value = 0.11
$ws1.Range("AB33").Value = 0.22
$ws1.Range("AC33").Value = 0.14
$ws1.Range("AB34").Value = 0.23
$ws1.Range("AC34").Value = 0.21
$ws1.Range("AB35").Value = 0.08
$ws1.Range("AC35").Value = 0.08
$ws1.Range("AB36").Value = 0.13
$ws1.Range("AC36").Value = 0.12
$ws1.Range("AB37").Value = 0.13
$ws1.Range("AC37").Value = 0.11
$ws1.Range("AB38").Value = 0.07
$ws1.Range("AC38").Value = 0.07
$ws1.Range("AB39").Value = 0.25
$ws1.Range("AC39").Value = 0.24
$ws1.Range("AB40").Value = 0.14
$ws1.Range("AC40").Value = 0.14
$ws1.Range("AB41").Value = 0.06
$ws1.Range("AC41").Value = 0.06
$ws1.Range("AB42").Value = 0.08
$ws1.Range("AC42").Value = 0.04
$ws1.Range("AB43").Value = 0.08
$ws1.Range("AC43").Value = 0.07
$ws1.Range("AB44").Value = 0.13
$ws1.Range("AC44").Value = 0.14
$ws1.Range("AB45").Value = 0.27
$ws1.Range("AC45").Value = 0.23

# --- Sheet "pocetR": add columns AA/AB for the two new survey waves ---
$ws2 = $wb.Worksheets.Item("pocetR")
$ws2.Range("AA1").Value = "13. 4. 2021"
$ws2.Range("AB1").Value = "4. 5. 2021"
$ws2.Range("Z1").Copy()
$ws2.Range("AA1:AB1").PasteSpecial(-4122)

$ws2.Range("AA2").Value = 2059
$ws2.Range("AB2").Value = 2032
$ws2.Range("AA3").Value = 231
$ws2.Range("AB3").Value = 218
$ws2.Range("AA4").Value = 455
$ws2.Range("AB4").Value = 427
$ws2.Range("AA5").Value = 1373
$ws2.Range("AB5").Value = 1387
$ws2.Range("AA6").Value = 983
$ws2.Range("AB6").Value = 966
$ws2.Range("AA7").Value = 183
$ws2.Range("AB7").Value = 181
$ws2.Range("AA8").Value = 589
$ws2.Range("AB8").Value = 581
$ws2.Range("AA9").Value = 304
$ws2.Range("AB9").Value = 304
$ws2.Range("AA10").Value = 936
$ws2.Range("AB10").Value = 915
$ws2.Range("AA11").Value = 169
$ws2.Range("AB11").Value = 168
$ws2.Range("AA12").Value = 129
$ws2.Range("AB12").Value = 133
$ws2.Range("AA13").Value = 825
$ws2.Range("AB13").Value = 816
$ws2.Range("AA14").Value = 953
$ws2.Range("AB14").Value = 936
$ws2.Range("AA15").Value = 702
$ws2.Range("AB15").Value = 697
$ws2.Range("AA16").Value = 404
$ws2.Range("AB16").Value = 399
$ws2.Range("AA17").Value = 246
$ws2.Range("AB17").Value = 240
$ws2.Range("AA18").Value = 781
$ws2.Range("AB18").Value = 736
$ws2.Range("AA19").Value = 645
$ws2.Range("AB19").Value = 682
$ws2.Range("AA20").Value = 252
$ws2.Range("AB20").Value = 254
$ws2.Range("AA21").Value = 520
$ws2.Range("AB21").Value = 576
$ws2.Range("AA22").Value = 372
$ws2.Range("AB22").Value = 340
$ws2.Range("AA23").Value = 234
$ws2.Range("AB23").Value = 190

# footer row 24 gets blank string cells in the new columns, matching the rest of the row
$ws2.Range("AA24").Value = ""
$ws2.Range("AB24").Value = ""

# --- Update the two "aktualizace" note strings (7. 4. 2021 -> 11. 5. 2021) ---
$ws1.Range("A46").Value = "Život během pandemie, Strategie domácností, % respondentů celkově a ve skupinách, aktualizace 11. 5. 2021"
$ws2.Range("A24").Value = "Život během pandemie, Strategie domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 11. 5. 2021"

Write-Output "edit applied"
